$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Bump the internal sheetId counter up to 6 so the new sheet we keep ---
# --- ends up with sheetId="6" (matching the target workbook.xml).       ---
$null = $wb.Worksheets.Add()
$null = $wb.Worksheets.Add()

# Remove the throwaway "Sheet1" (sheetId 5), and the old "testdata"/"system" sheets.
$wb.Worksheets("Sheet1").Delete()
$wb.Worksheets("testdata").Delete()
$wb.Worksheets("system").Delete()

# "Sheet2" (sheetId 6) is our new sheet; move it right after "hub" and rename it.
$wb.Worksheets("Sheet2").Move($null, $wb.Worksheets("hub"))
$scr = $wb.Worksheets("Sheet2")
$scr.Name = "screenTitles"

# --- Populate the screenTitles sheet with its header + data rows ---
$scr.Range("A1").Value2 = "objectID"
$scr.Range("B1").Value2 = "name_nl"

$scr.Range("B4").Value2 = "epg"
$scr.Range("B2").Value2 = "instellingen"
$scr.Range("A3").Value2 = "System"
$scr.Range("A2").Value2 = "Setting"
$scr.Range("A4").Value2 = "epgSetting"
$scr.Range("B3").Value2 = "systeem"

# Header row style (bold/shaded header, same style used elsewhere in the workbook)
$wb.Worksheets("hub").Range("A1:B1").Copy()
$scr.Range("A1:B1").PasteSpecial(-4122)

# Column widths for the new sheet
$scr.Columns("A").ColumnWidth = 14.67
$scr.Columns("B").ColumnWidth = 17.83

# hub sheet selection moves to A7 (and loses being the tab-selected sheet)
$wb.Worksheets("hub").Range("A7").Select()

# View: screenTitles becomes the selected/active sheet (must be done last so
# it ends up as the active tab in the saved workbook)
$scr.Range("T25").Select()
$scr.Activate()
